# Updated cryptos list on Sat Jun  1 10:00:12 UTC 2024 with GitHub Actions
#
# Refreshes the per-row crypto Price (column D) and Volume(1h) (column E)
# figures pulled from coinranking.com, and restores the correct row order
# for RenderToken / Dai (rows 28-29 had swapped on the previous run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.711.50"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "'3.792.72"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'595.60"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "'166.99"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'3.792.03"
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("D11").Value = "'6.32"
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "'4.426.24"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").Value = "'3.803.60"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").Value = "'67.681.06"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "'18.41"
$ws.Range("E18").Value = "  +3.36%  "
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("E21").Value = "  -6.30%  "
$ws.Range("D22").Value = "'459.14"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").Value = "'0.697"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("E24").Value = "  +4.68%  "
$ws.Range("D25").Value = "'83.30"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("E26").Value = "  +2.07%  "
$ws.Range("D27").Value = "'2.12"
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.01"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("E31").Value = "  +5.13%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "'29.62"
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").Value = "'3.736.94"
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("D37").Value = "'0.1000"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("D38").Value = "'3.38"
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D44").Value = "'44.11"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "'48.03"
$ws.Range("E45").Value = "  +2.73%  "
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "'149.65"
$ws.Range("E47").Value = "  +3.25%  "
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").Value = "'391.51"
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("D50").Value = "'26.76"
$ws.Range("E50").Value = "  +7.67%  "
$ws.Range("E51").Value = "  -4.27%  "
